$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in column H, copying the formatting used by the
# other header cells (bold, bordered, centered style already on G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Label"

# Fill in Label values: 0 for Control-diagnosis rows, 1 for MDD-diagnosis rows
$labels = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 1
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 8).Value = $labels[$row]
}
